# Auto commit at 2025-09-25  8:03:46.41
# Append two new daily rows (2025-09-24 / serial 45924) for the two
# charging stations, continuing the existing daily log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Date serial 45924 == 2025-09-24 (row 47 above it is 45923 == 2025-09-23).
$dateSerial = 45924

# Row 48 - 四方坪站
$ws.Cells.Item(48, 1).Value = $dateSerial
$ws.Cells.Item(48, 2).Value = "四方坪站"
$ws.Cells.Item(48, 3).Value = 8012.71
$ws.Cells.Item(48, 4).Value = 6438.57
$ws.Cells.Item(48, 5).Value = 2768.25
$ws.Cells.Item(48, 6).Value = 361

# Row 49 - 高岭站
$ws.Cells.Item(49, 1).Value = $dateSerial
$ws.Cells.Item(49, 2).Value = "高岭站"
$ws.Cells.Item(49, 3).Value = 4060.05
$ws.Cells.Item(49, 4).Value = 3156.64
$ws.Cells.Item(49, 5).Value = 1052.47
$ws.Cells.Item(49, 6).Value = 152

# Give the new rows the same formatting as the row directly above them
# (row 47) by copying just the cell formats, not the values/dates.
$ws.Range("A47:F47").Copy()
$ws.Range("A48:F49").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Mirror the scrolled/selected state shown after the new rows were added.
$ws.Range("H47").Select()
